$wb = $excel.ActiveWorkbook

# Add a new worksheet that will hold the "discount_rate" data
$ws = $wb.Worksheets.Add()
$ws.Name = "discount_rate"

# Header row (mirrors the layout used by the OPEX sheet)
$ws.Range("B1").Value = "discount_rate"
$ws.Range("C1").Value = "unit"
$ws.Range("D1").Value = "reference"
$ws.Range("E1").Value = "comment"

# Technology rows, same technologies/order as OPEX, but with a 7% discount rate
$ws.Range("A2").Value = "Photovoltaic"
$ws.Range("B2").Value = 0.07
$ws.Range("C2").Value = "fraction of capex"

$ws.Range("A3").Value = "WindTurbine_Onshore_4000"
$ws.Range("B3").Value = 0.07
$ws.Range("C3").Value = "fraction of capex"

$ws.Range("A4").Value = "NuclearPlant"
$ws.Range("B4").Value = 0.07
$ws.Range("C4").Value = "fraction of capex"

$ws.Range("A5").Value = "Storage_Battery"
$ws.Range("B5").Value = 0.07
$ws.Range("C5").Value = "fraction of capex"

$ws.Range("A6").Value = "Hydro_Reservoir"
$ws.Range("B6").Value = 0.07
$ws.Range("C6").Value = "fraction of capex"

$ws.Range("A7").Value = "GasTurbine_simple"
$ws.Range("B7").Value = 0.07
$ws.Range("C7").Value = "fraction of capex"

$ws.Range("A8").Value = "CoalPlant"
$ws.Range("B8").Value = 0.07
$ws.Range("C8").Value = "fraction of capex"

$ws.Range("A9").Value = "PumpedHydro_Closed"
$ws.Range("B9").Value = 0.07
$ws.Range("C9").Value = "fraction of capex"

# Column widths, matching the OPEX sheet styling
$ws.Columns.Item(1).ColumnWidth = 32.5546875
$ws.Columns.Item(2).ColumnWidth = 8.88671875
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 19.6640625
$ws.Columns.Item(5).ColumnWidth = 35.77734375

# Selection matching the final diff (B2:B9 selected, B2 active)
$ws.Range("B2:B9").Select()

# Move the new sheet to the very first position in the workbook
$ws.Move($wb.Worksheets.Item(1))

# Make it the active/selected sheet, as it is the first tab in the final workbook
$ws.Select()
